# Weekly driver report update for 2025-04-20
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Bad Drivers" section ---
$ws.Range("D3").Value = 80.09999999999999
$ws.Range("C4").Value = 4970
$ws.Range("D4").Value = 90.7
$ws.Range("C5").Value = 6466

# --- "Good Drivers" section (rows 13-26) : weekly refresh / re-sort ---
# Column E holds the "Driver Vintage" as literal text (e.g. "2024-11-10"),
# not a real date. Excel's smart Range.Value parser auto-converts date-like
# strings into date serials (and swaps in a date NumberFormat) unless the
# cell is pre-formatted as Text ("@"). We set "@" first, assign the literal
# text, then restore the original (General / right-aligned) look by pasting
# just the formatting from a stable same-style reference cell (B3, whose own
# value never changes in this edit) so the cell's style index is preserved.
$ws.Range("B3").Copy() | Out-Null

$ws.Range("A13").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.0.1.5"
$ws.Range("B13").Value = 156943
$ws.Range("D13").Value = 100

$ws.Range("A14").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.120.0.3"
$ws.Range("B14").Value = 34181
$ws.Range("D14").Value = 99.90000000000001
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "2025-02-05"
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("A15").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 23.100.0.4"
$ws.Range("B15").Value = 445055
$ws.Range("D15").Value = 99.90000000000001
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2024-11-10"
$ws.Range("E15").PasteSpecial(-4122)

$ws.Range("A16").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 23.20.1.1"
$ws.Range("B16").Value = 13533
$ws.Range("D16").Value = 100
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2023-12-19"
$ws.Range("E16").PasteSpecial(-4122)

$ws.Range("A17").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.170.2.1"
$ws.Range("B17").Value = 19083
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2022-08-30"
$ws.Range("E17").PasteSpecial(-4122)

$ws.Range("A18").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.100.0.3"
$ws.Range("B18").Value = 12988
$ws.Range("D18").Value = 100
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "2022-05-01"
$ws.Range("E18").PasteSpecial(-4122)

$ws.Range("A19").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.130.0.5"
$ws.Range("B19").Value = 18738
$ws.Range("D19").Value = 99.90000000000001
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "2022-03-14"
$ws.Range("E19").PasteSpecial(-4122)

$ws.Range("A20").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.110.1.1"
$ws.Range("B20").Value = 42024
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2022-01-01"
$ws.Range("E20").PasteSpecial(-4122)

$ws.Range("A21").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.80.0.9"
$ws.Range("B21").Value = 77849
$ws.Range("D21").Value = 99.90000000000001
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2021-08-18"
$ws.Range("E21").PasteSpecial(-4122)

$ws.Range("A22").Value = "Intel(R) Wi-Fi 6E AX210 160MHz - 22.70.0.6"
$ws.Range("B22").Value = 15504
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2021-06-28"
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("A23").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 22.50.1.1"
$ws.Range("B23").Value = 34244
$ws.Range("D23").Value = 100
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "2021-04-27"
$ws.Range("E23").PasteSpecial(-4122)

$ws.Range("A24").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.110.3.2"
$ws.Range("B24").Value = 59673
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2020-08-05"
$ws.Range("E24").PasteSpecial(-4122)

$ws.Range("A25").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.70.0.6"
$ws.Range("B25").Value = 113652
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2020-01-06"
$ws.Range("E25").PasteSpecial(-4122)

$ws.Range("A26").Value = "Intel(R) Wi-Fi 6 AX201 160MHz - 21.60.2.1"
$ws.Range("B26").Value = 56018
